$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
$ws.Range("A1").Value = "Country"
$ws.Range("B1").Value = "Subregion"
$ws.Range("C1").Value = "Poplation"
$ws.Range("D1").Value = "European Union Association"

# Copy the header style (bold/centered/bordered) from C1 to the new D1 header cell
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Fill in column D (European Union Association) for each data row ---
$associations = @{
    2  = "No Association"
    3  = "European Union (EU)"
    4  = "No Association"
    5  = "European Union (EU)"
    6  = "No Association"
    7  = "European Union (EU)"
    8  = "European Union (EU)"
    9  = "European Union (EU)"
    10 = "No Association"
    11 = "European Union (EU)"
    12 = "European Union (EU)"
    13 = "European Union (EU)"
    14 = "European Union (EU)"
    15 = "European Union (EU)"
    16 = "European Union (EU)"
    17 = "European Union (EU)"
    18 = "European Free Trade Association (EFTA)"
    19 = "European Union (EU)"
    20 = "European Union (EU)"
    21 = "No Association"
    22 = "No Association"
    23 = "European Union (EU)"
    24 = "European Union (EU)"
    25 = "European Union (EU)"
    26 = "European Union (EU)"
    27 = "No Association"
    28 = "No Association"
    29 = "European Union (EU)"
    30 = "No Association"
    31 = "European Free Trade Association (EFTA)"
    32 = "European Union (EU)"
    33 = "European Union (EU)"
    34 = "European Union (EU)"
    35 = "No Association"
    36 = "No Association"
    37 = "European Union (EU)"
    38 = "European Union (EU)"
    39 = "European Union (EU)"
    40 = "European Union (EU)"
    41 = "European Free Trade Association (EFTA)"
    42 = "No Association"
    43 = "No Association"
}

foreach ($row in $associations.Keys) {
    $ws.Cells.Item($row, 4).Value = $associations[$row]
}
